$wb = $excel.ActiveWorkbook

# --- Sheet1: bump the numeric suffix used by the CONCATENATE formulas ---
$calc = $wb.Worksheets.Item("Sheet1")
$calc.Range("I2").Value = 14
$calc.Range("I23").Value = 19

# --- login sheet: usernames/passwords/emails ending in "13" -> "14" ---
$login = $wb.Worksheets.Item("login")
for ($r = 2; $r -le 21; $r++) {
    $gCell = $login.Cells.Item($r, 7)
    $hCell = $login.Cells.Item($r, 8)
    $iCell = $login.Cells.Item($r, 9)
    $gVal = $gCell.Value()
    $hVal = $hCell.Value()
    $iVal = $iCell.Value()
    $gCell.Value = $gVal.Replace("13", "14")
    $hCell.Value = $hVal.Replace("13", "14")
    $iCell.Value = $iVal.Replace("13", "14")
}

# --- order sheet: usernames/passwords/emails ending in "18" -> "19" ---
$order = $wb.Worksheets.Item("order")
for ($r = 2; $r -le 21; $r++) {
    $rCell = $order.Cells.Item($r, 18)
    $sCell = $order.Cells.Item($r, 19)
    $tCell = $order.Cells.Item($r, 20)
    $rVal = $rCell.Value()
    $sVal = $sCell.Value()
    $tVal = $tCell.Value()
    $rCell.Value = $rVal.Replace("18", "19")
    $sCell.Value = $sVal.Replace("18", "19")
    $tCell.Value = $tVal.Replace("18", "19")
}
